# Fruta / hortaliza, semanal
# Insert a new weekly price record for Naranja (Valencia, Tercera) above the
# existing row 126, shifting the subsequent rows (old 126-136) down to 127-137.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 126; Excel shifts rows 126:136 -> 127:137
# and extends the used range to A1:T137 automatically.
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new record.
$ws.Cells.Item(126, 1).Value = 1
$ws.Cells.Item(126, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(126, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(126, 4).Value = 45021
$ws.Cells.Item(126, 5).Value = 15
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100102
$ws.Cells.Item(126, 8).Value = "Cítricos"
$ws.Cells.Item(126, 9).Value = 100102005
$ws.Cells.Item(126, 10).Value = "Naranja"
$ws.Cells.Item(126, 11).Value = "Valencia"
$ws.Cells.Item(126, 12).Value = "Tercera"
$ws.Cells.Item(126, 13).Value = 270
$ws.Cells.Item(126, 14).Value = 1300
$ws.Cells.Item(126, 15).Value = 1400
$ws.Cells.Item(126, 16).Value = 1350
$ws.Cells.Item(126, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(126, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(126, 19).Value = 1350
$ws.Cells.Item(126, 20).Value = 1
